# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计", built by copying
#    the existing "2022-Q1" sheet (so header styling / layout matches) and
#    then overwriting its values with the 2022-Q4 numbers (6 rows incl.
#    header; 2 more data rows than the source sheet, so their formatting
#    is extended explicitly before being filled in).
# 2. Update the "总计" summary sheet: existing row 2 (2022-Q1 totals)
#    becomes the new 2022-Q4 totals, and the old rows shift down one:
#    row 3 = 2022-Q1, row 4 = 2021-Q4 (previously rows 2 and 3).

$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q4" sheet -------------------------------

$source = $wb.Worksheets.Item("2022-Q1")
$anchor = $wb.Worksheets.Item("总计")
$source.Copy($null, $anchor)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Extend the row-style (bordered "index" column) from row 2 down to the two
# new rows (5 and 6) before writing values into them.
$q4.Range("A2").Copy()
$q4.Range("A5:A6").PasteSpecial(-4122)

# Fund-code column (B) and the numeric-but-text data columns (D:G) need to
# stay text (leading zeros / fixed decimal strings) - format the full data
# range as text up front so assigning the literal strings below doesn't get
# auto-coerced into numbers.
$q4.Range("B2:B6").NumberFormat = "@"
$q4.Range("D2:G6").NumberFormat = "@"

# Header row (row 1) is already correct (copied verbatim from 2022-Q1).

# Row 2: 013442 / 建信中证1000指数增强E
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "013442"
$q4.Range("C2").Value = "建信中证1000指数增强E"
$q4.Range("D2").Value = "9.52"
$q4.Range("E2").Value = "86.80"
$q4.Range("F2").Value = "1.38"
$q4.Range("G2").Value = "0.1314"
$q4.Range("H2").Value = 3

# Row 3: 006165 / 建信中证1000指数增强A
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "006165"
$q4.Range("C3").Value = "建信中证1000指数增强A"
$q4.Range("D3").Value = "7.20"
$q4.Range("E3").Value = "86.80"
$q4.Range("F3").Value = "1.38"
$q4.Range("G3").Value = "0.0994"
$q4.Range("H3").Value = 3

# Row 4: 006166 / 建信中证1000指数增强C
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "006166"
$q4.Range("C4").Value = "建信中证1000指数增强C"
$q4.Range("D4").Value = "2.21"
$q4.Range("E4").Value = "86.80"
$q4.Range("F4").Value = "1.38"
$q4.Range("G4").Value = "0.0305"
$q4.Range("H4").Value = 3

# Row 5: 003241 / 创金合信量化发现灵活配置混合A
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "003241"
$q4.Range("C5").Value = "创金合信量化发现灵活配置混合A"
$q4.Range("D5").Value = "0.30"
$q4.Range("E5").Value = "91.60"
$q4.Range("F5").Value = "0.93"
$q4.Range("G5").Value = "0.0028"
$q4.Range("H5").Value = 7

# Row 6: 003242 / 创金合信量化发现灵活配置混合C
$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "003242"
$q4.Range("C6").Value = "创金合信量化发现灵活配置混合C"
$q4.Range("D6").Value = "0.29"
$q4.Range("E6").Value = "91.60"
$q4.Range("F6").Value = "0.93"
$q4.Range("G6").Value = "0.0027"
$q4.Range("H6").Value = 7

# --- Step 2: update the "总计" summary sheet -------------------------------

$total = $wb.Worksheets.Item("总计")

# Row 4 (new) needs the bordered "index" style (s="2") applied, same as
# rows 2 and 3 - copy that formatting down before writing the 2021-Q4
# values there.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)

# Write rows bottom-up so no in-progress values are overwritten by a later
# statement referencing the same cell.
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.05

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.04

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.27
